$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.975.98'
$ws.Range('E2').Value = '  +1.68%  '
$ws.Range('D3').Value = '1.905.09'
$ws.Range('E3').Value = '  +0.81%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.96'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('E6').Value = '  -0.80%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.43'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.364'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '57.40'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +8.09%  '
$ws.Range('E11').Value = '  +2.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0992'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.71'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +12.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.794'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.23%  '
$ws.Range('D15').Value = '2.181.42'
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.10'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +3.15%  '
$ws.Range('D17').Value = '1.895.42'
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('D18').Value = '35.977.78'
$ws.Range('E18').Value = '  +1.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.91'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.86%  '
$ws.Range('D20').Value = '0.0₃0841'
$ws.Range('E20').Value = '  +1.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '249.74'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.11'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.24'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +5.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.71'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.49%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.24'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.36'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.65%  '
$ws.Range('E28').Value = '  +3.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.70'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.71%  '
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.53'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +6.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0606'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.32'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.25%  '
$ws.Range('E34').Value = '  +3.17%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.47'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -16.60%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0830'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +19.74%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('E39').Value = '  +0.52%  '
$ws.Range('E40').Value = '  +4.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '102.17'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +5.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '15.19'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +21.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.83'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.47%  '
$ws.Range('E44').Value = '  +1.92%  '
$ws.Range('D45').Value = '1.332.72'
$ws.Range('E45').Value = '  +2.96%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.36'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0807'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.22%  '
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.78'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.39'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.81%  '
$ws.Range('D51').Value = '2.078.60'
$ws.Range('E51').Value = '  +0.38%  '
